$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert "Common Dab" row before the current "Flatfish" row (row 5),
# and "Plaice" row before the current "Scorpionfish" row (now row 7 after
# the first insert).
$ws.Rows("5").Insert()
$ws.Range("A5").Value = "Common Dab"
$ws.Range("B5").Value = 1.5
$ws.Range("C5").Value = 0.3
$ws.Range("D5").Value = 1.7
$ws.Range("E5").Value = 0.2

$ws.Rows("7").Insert()
$ws.Range("A7").Value = "Plaice"
$ws.Range("B7").Value = 1.5
$ws.Range("C7").Value = 0.3
$ws.Range("D7").Value = 1.7
$ws.Range("E7").Value = 0.2

# Column A width widened slightly (author manually dragged the column
# border wider after the new, longer species names were added)
$ws.Columns("A").ColumnWidth = 11.67

# Selection moved
$ws.Range("E8").Select()
